$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 74.75
$ws.Range("I42").Value = 100
$ws.Range("J42").Value = 66.333336
$ws.Range("K42").Value = 300
$ws.Range("L42").Value = 199.000008
$ws.Range("M42").Value = -70
$ws.Range("N42").Value = -659.000008
# Row 108
$ws.Range("H108").Value = 15550
$ws.Range("J108").Value = 15550
$ws.Range("L108").Value = 15550
$ws.Range("N108").Value = -23230
# Row 132
$ws.Range("H132").Value = 4338.4614
$ws.Range("I132").Value = 2939.4
$ws.Range("J132").Value = 9002
$ws.Range("K132").Value = 8818.200000000001
$ws.Range("L132").Value = 27006
$ws.Range("M132").Value = -6288.200000000001
$ws.Range("N132").Value = -32066
# Row 137
$ws.Range("H137").Value = 1920.7142
$ws.Range("I137").Value = 1971.625
$ws.Range("K137").Value = 5914.875
$ws.Range("M137").Value = -3364.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 172.5
$ws.Range("I5").Value = 213.5
$ws.Range("K5").Value = 213.5
$ws.Range("M5").Value = -101.5
# Row 122
$ws.Range("H122").Value = 8639.857
$ws.Range("J122").Value = 5248.75
$ws.Range("L122").Value = 15746.25
$ws.Range("N122").Value = -20646.25
# Row 132
$ws.Range("H132").Value = 1628.5385
$ws.Range("I132").Value = 1523.091
$ws.Range("J132").Value = 2208.5
$ws.Range("K132").Value = 4569.272999999999
$ws.Range("L132").Value = 6625.5
$ws.Range("M132").Value = -2039.272999999999
$ws.Range("N132").Value = -11685.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 172.5
$ws.Range("I4").Value = 213.5
$ws.Range("K4").Value = 213.5
$ws.Range("M4").Value = -98.5
# Row 80
$ws.Range("H80").Value = 935.3333
$ws.Range("J80").Value = 1251.5
$ws.Range("L80").Value = 1251.5
$ws.Range("N80").Value = -3247.5
# Row 83
$ws.Range("H83").Value = 935.3333
$ws.Range("J83").Value = 1251.5
$ws.Range("L83").Value = 6257.5
$ws.Range("N83").Value = -16241.5
# Row 99
$ws.Range("H99").Value = 3343.2354
$ws.Range("I99").Value = 3427.1875
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 3427.1875
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -1929.1875
$ws.Range("N99").Value = -4996
# Row 134
$ws.Range("H134").Value = 5979.9287
$ws.Range("J134").Value = 9365
$ws.Range("L134").Value = 28095
$ws.Range("N134").Value = -33165

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10378.333
$ws.Range("I31").Value = 11235.333
$ws.Range("J31").Value = 8664.333000000001
$ws.Range("K31").Value = 11235.333
$ws.Range("L31").Value = 8664.333000000001
$ws.Range("M31").Value = -10940.333
$ws.Range("N31").Value = -9254.333000000001
# Row 34
$ws.Range("H34").Value = 10378.333
$ws.Range("I34").Value = 11235.333
$ws.Range("J34").Value = 8664.333000000001
$ws.Range("K34").Value = 11235.333
$ws.Range("L34").Value = 8664.333000000001
$ws.Range("M34").Value = -11033.333
$ws.Range("N34").Value = -9068.333000000001
# Row 105
$ws.Range("H105").Value = 1151.6666
$ws.Range("J105").Value = 850
$ws.Range("L105").Value = 850
$ws.Range("N105").Value = -4344

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 629.6667
$ws.Range("I11").Value = 629.6667
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1889.0001
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1749.0001
$ws.Range("N11").ClearContents()
# Row 108
$ws.Range("H108").Value = 167.5
$ws.Range("I108").Value = 167.5
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 502.5
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 2377.5
$ws.Range("N108").ClearContents()
# Row 128
$ws.Range("H128").Value = 125000
$ws.Range("I128").Value = 125000
$ws.Range("K128").Value = 375000
$ws.Range("M128").Value = -370020

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3485.625
$ws.Range("I102").Value = 3312.4285
$ws.Range("K102").Value = 3312.4285
$ws.Range("M102").Value = -1690.4285
# Row 132
$ws.Range("H132").Value = 2647.3684
$ws.Range("I132").Value = 2147.5
$ws.Range("J132").Value = 5313.3335
$ws.Range("K132").Value = 6442.5
$ws.Range("L132").Value = 15940.0005
$ws.Range("M132").Value = -3912.5
$ws.Range("N132").Value = -21000.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3792.375
$ws.Range("I40").Value = 3792.375
$ws.Range("K40").Value = 3792.375
$ws.Range("M40").Value = -3656.375
# Row 88
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30856
# Row 91
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32964
# Row 93
$ws.Range("H93").Value = 4477.5713
$ws.Range("I93").Value = 4688.6
$ws.Range("J93").Value = 3950
$ws.Range("K93").Value = 4688.6
$ws.Range("L93").Value = 3950
$ws.Range("M93").Value = -3440.6
$ws.Range("N93").Value = -6446
# Row 100
$ws.Range("H100").Value = 7054.4165
$ws.Range("I100").Value = 3385.5715
$ws.Range("J100").Value = 12190.8
$ws.Range("K100").Value = 3385.5715
$ws.Range("L100").Value = 12190.8
$ws.Range("M100").Value = -2844.5715
$ws.Range("N100").Value = -13272.8
# Row 122
$ws.Range("H122").Value = 3271.2856
$ws.Range("I122").Value = 2649.8333
$ws.Range("K122").Value = 7949.499899999999
$ws.Range("M122").Value = -5499.499899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 5001000
$ws.Range("J5").Value = 5001000
$ws.Range("L5").Value = 5001000
$ws.Range("N5").Value = -5001224
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 122
$ws.Range("H122").Value = 1719.2
$ws.Range("I122").Value = 1098.1428
$ws.Range("K122").Value = 3294.4284
$ws.Range("M122").Value = -844.4284000000002
# Row 126
$ws.Range("H126").Value = 5497.25
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5
# Row 136
$ws.Range("H136").Value = 1845.625
$ws.Range("I136").Value = 1909.2858
$ws.Range("K136").Value = 5727.857400000001
$ws.Range("M136").Value = -3177.857400000001
